# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.418.12"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = "'3.507.81"
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'586.77"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'135.90"
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("D7").Value = "'3.508.80"
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = "'0.375"
$ws.Range("E12").Value = '  -3.62%  '
$ws.Range("D13").Value = "'4.104.43"
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  +0.05%  '
$ws.Range("D15").Value = "'0.119"
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").Value = "'3.507.55"
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = "'64.383.90"
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = "'25.19"
$ws.Range("E18").Value = '  -9.56%  '
$ws.Range("D19").Value = "'9.86"
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").Value = "'13.86"
$ws.Range("E20").Value = '  -2.85%  '
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = "'384.38"
$ws.Range("E22").Value = '  -1.79%  '
$ws.Range("D23").Value = "'0.571"
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").Value = "'3.646.46"
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").Value = "'74.02"
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +1.20%  '
$ws.Range("E28").Value = '  +3.82%  '
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("E30").Value = '  +1.04%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = "'8.29"
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = "'3.527.60"
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D36").Value = "'0.148"
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").Value = "'23.59"
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D38").Value = "'5.31"
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").Value = "'163.88"
$ws.Range("E41").Value = '  -4.31%  '
$ws.Range("D42").Value = "'0.0787"
$ws.Range("E42").Value = '  -2.97%  '
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("D44").Value = "'25.95"
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = "'41.94"
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").Value = "'4.42"
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").Value = "'2.483.33"
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("E51").Value = '  -1.89%  '
